$d = $word.ActiveDocument

# --- Paragraph 42: "HIVE TEAM: QUALITY ASSURANCE" heading ---
# This paragraph starts with a standalone run containing a single
# (non-breaking) space in Tahoma, followed by a run with the title text
# in Helvetica Neue. The edit drops the leading space-only run entirely
# and retranslates the title run into Russian.
$pTitle = $d.Paragraphs(42)
$rTitle = $pTitle.Range
$leadChar = $d.Range($rTitle.Start, $rTitle.Start + 1)
$leadChar.Text = ""

$pTitle = $d.Paragraphs(42)
$pTitle.Range.Find.Execute("HIVE TEAM: QUALITY ASSURANCE", $true, $false, $false, $false, $false, $true, 1, $false, "Команда Hive: Контроль Качества", 2)

# --- Paragraph 43: section description ---
$d.Paragraphs(43).Range.Find.Execute("Ensuring all development tasks meet quality criteria.", $true, $false, $false, $false, $false, $true, 1, $false, "Тестирование и контроль качества на всех этапах разработки.", 2)

# --- Paragraph 45: "Hive Coordinator" role (Cryptolize) ---
$d.Paragraphs(45).Range.Find.Execute("Hive Coordinator", $true, $false, $false, $false, $false, $true, 1, $false, "Координатор Hive", 2)

# --- Paragraph 47: "Release Coordinator" role (Jazz) ---
$d.Paragraphs(47).Range.Find.Execute("Release Coordinator", $true, $false, $false, $false, $false, $true, 1, $false, "Релиз-координатор", 2)

# --- Paragraph 49: "auditor" role (emelia) ---
$d.Paragraphs(49).Range.Find.Execute("auditor", $true, $false, $false, $false, $false, $true, 1, $false, "Аудитор", 2)

# --- Paragraph 51: "Senior QA Tester" role (Nitego) ---
$d.Paragraphs(51).Range.Find.Execute("Senior QA Tester", $true, $false, $false, $false, $false, $true, 1, $false, "Специалист по обеспечению качества", 2)
